# Rename "SwateTemplateMetadata" sheet to "isa_template" and make it the
# active/selected sheet (tab) in the workbook.

$wb = $excel.ActiveWorkbook

$metaSheet = $wb.Worksheets.Item("SwateTemplateMetadata")
$metaSheet.Name = "isa_template"

$metaSheet.Activate()
$metaSheet.Select()
